$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.066553863589701
$ws.Cells.Item(2, 4).Value = 1.067738104166037
$ws.Cells.Item(2, 5).Value = 1.071153619569611
$ws.Cells.Item(2, 6).Value = 1.080852740821386
$ws.Cells.Item(2, 9).Value = 1.052233023249721
$ws.Cells.Item(2, 10).Value = 1.071503274683289
$ws.Cells.Item(2, 11).Value = 1.070445162377491
$ws.Cells.Item(2, 12).Value = 1.07385155035598
$ws.Cells.Item(2, 13).Value = 1.083525091945799
$ws.Cells.Item(2, 14).Value = 1.073024932236483

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.067706599063331
$ws.Cells.Item(3, 4).Value = 1.068640652648386
$ws.Cells.Item(3, 5).Value = 1.072172123426892
$ws.Cells.Item(3, 6).Value = 1.08193232469576
$ws.Cells.Item(3, 9).Value = 1.052565551092424
$ws.Cells.Item(3, 10).Value = 1.072311083910006
$ws.Cells.Item(3, 11).Value = 1.071163338258711
$ws.Cells.Item(3, 12).Value = 1.074686053911982
$ws.Cells.Item(3, 13).Value = 1.084422377585229
$ws.Cells.Item(3, 14).Value = 1.073833888644959

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.068452645660686
$ws.Cells.Item(4, 4).Value = 1.069224768824274
$ws.Cells.Item(4, 5).Value = 1.072831599501723
$ws.Cells.Item(4, 6).Value = 1.082631444441376
$ws.Cells.Item(4, 9).Value = 1.052779594632019
$ws.Cells.Item(4, 10).Value = 1.072833354482821
$ws.Cells.Item(4, 11).Value = 1.071627516663699
$ws.Cells.Item(4, 12).Value = 1.075225855500461
$ws.Cells.Item(4, 13).Value = 1.085002934186662
$ws.Cells.Item(4, 14).Value = 1.074356900901892

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.068766320148599
$ws.Cells.Item(5, 4).Value = 1.069470356717789
$ws.Cells.Item(5, 5).Value = 1.073108947586359
$ws.Cells.Item(5, 6).Value = 1.08292548765906
$ws.Cells.Item(5, 9).Value = 1.052869309412432
$ws.Cells.Item(5, 10).Value = 1.073052812792247
$ws.Cells.Item(5, 11).Value = 1.071822530541113
$ws.Cells.Item(5, 12).Value = 1.075452745519223
$ws.Cells.Item(5, 13).Value = 1.085246988896202
$ws.Cells.Item(5, 14).Value = 1.07457667086729

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.068818989635191
$ws.Cells.Item(6, 4).Value = 1.069511593493945
$ws.Cells.Item(6, 5).Value = 1.073155521668049
$ws.Cells.Item(6, 6).Value = 1.082974866624056
$ws.Cells.Item(6, 9).Value = 1.052884357137763
$ws.Cells.Item(6, 10).Value = 1.073089654699642
$ws.Cells.Item(6, 11).Value = 1.071855266798665
$ws.Cells.Item(6, 12).Value = 1.075490838856918
$ws.Cells.Item(6, 13).Value = 1.085287966094544
$ws.Cells.Item(6, 14).Value = 1.074613565094419

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.06845683685028
$ws.Cells.Item(7, 4).Value = 1.069228050281901
$ws.Cells.Item(7, 5).Value = 1.07283530503044
$ws.Cells.Item(7, 6).Value = 1.082635372937574
$ws.Cells.Item(7, 9).Value = 1.052780794462754
$ws.Cells.Item(7, 10).Value = 1.072836287304896
$ws.Cells.Item(7, 11).Value = 1.071630122946585
$ws.Cells.Item(7, 12).Value = 1.075228887383199
$ws.Cells.Item(7, 13).Value = 1.085006195301956
$ws.Cells.Item(7, 14).Value = 1.074359837888911

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.066943405352074
$ws.Cells.Item(8, 4).Value = 1.068043102413865
$ws.Cells.Item(8, 5).Value = 1.071497737427243
$ws.Cells.Item(8, 6).Value = 1.081217475947299
$ws.Cells.Item(8, 9).Value = 1.052345635374524
$ws.Cells.Item(8, 10).Value = 1.071776367958851
$ws.Cells.Item(8, 11).Value = 1.070687982908504
$ws.Cells.Item(8, 12).Value = 1.074133611627391
$ws.Cells.Item(8, 13).Value = 1.083828343612515
$ws.Cells.Item(8, 14).Value = 1.07329841333582

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.064277661795934
$ws.Cells.Item(9, 4).Value = 1.06595589831709
$ws.Cells.Item(9, 5).Value = 1.069144112313488
$ws.Cells.Item(9, 6).Value = 1.07872322502271
$ws.Cells.Item(9, 9).Value = 1.051570218521204
$ws.Cells.Item(9, 10).Value = 1.069905308541703
$ws.Cells.Item(9, 11).Value = 1.069023755905249
$ws.Cells.Item(9, 12).Value = 1.072202229888017
$ws.Cells.Item(9, 13).Value = 1.081752455756646
$ws.Cells.Item(9, 14).Value = 1.071424696799687

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.062501205068692
$ws.Cells.Item(10, 4).Value = 1.064564989757302
$ws.Cells.Item(10, 5).Value = 1.067577272681671
$ws.Cells.Item(10, 6).Value = 1.077063255891244
$ws.Cells.Item(10, 9).Value = 1.05104748191494
$ws.Cells.Item(10, 10).Value = 1.068655673950294
$ws.Cells.Item(10, 11).Value = 1.067911540227503
$ws.Cells.Item(10, 12).Value = 1.070913717578419
$ws.Cells.Item(10, 13).Value = 1.080368282341687
$ws.Cells.Item(10, 14).Value = 1.070173287583823

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.061732136266216
$ws.Cells.Item(11, 4).Value = 1.063962844397286
$ws.Cells.Item(11, 5).Value = 1.066899344680956
$ws.Cells.Item(11, 6).Value = 1.076345148887813
$ws.Cells.Item(11, 9).Value = 1.050819756391215
$ws.Cells.Item(11, 10).Value = 1.068114028509633
$ws.Cells.Item(11, 11).Value = 1.067429289140665
$ws.Cells.Item(11, 12).Value = 1.070355556074504
$ws.Cells.Item(11, 13).Value = 1.079768859279033
$ws.Cells.Item(11, 14).Value = 1.069630872944508

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.061446491173447
$ws.Cells.Item(12, 4).Value = 1.063739199841733
$ws.Cells.Item(12, 5).Value = 1.066647610526468
$ws.Cells.Item(12, 6).Value = 1.076078512479895
$ws.Cells.Item(12, 9).Value = 1.050734961942458
$ws.Cells.Item(12, 10).Value = 1.067912754939413
$ws.Cells.Item(12, 11).Value = 1.067250061066991
$ws.Cells.Item(12, 12).Value = 1.070148195624924
$ws.Cells.Item(12, 13).Value = 1.079546196562466
$ws.Cells.Item(12, 14).Value = 1.069429313542732

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.061507762068052
$ws.Cells.Item(13, 4).Value = 1.063787171497566
$ws.Cells.Item(13, 5).Value = 1.066701604808842
$ws.Cells.Item(13, 6).Value = 1.076135702350224
$ws.Cells.Item(13, 9).Value = 1.050753160011424
$ws.Cells.Item(13, 10).Value = 1.06795593253569
$ws.Cells.Item(13, 11).Value = 1.067288510560498
$ws.Cells.Item(13, 12).Value = 1.070192676707171
$ws.Cells.Item(13, 13).Value = 1.079593958941354
$ws.Cells.Item(13, 14).Value = 1.069472552456148

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.061708524325127
$ws.Cells.Item(14, 4).Value = 1.063944357480529
$ws.Cells.Item(14, 5).Value = 1.066878534657319
$ws.Cells.Item(14, 6).Value = 1.076323106590233
$ws.Cells.Item(14, 9).Value = 1.050812751480254
$ws.Cells.Item(14, 10).Value = 1.068097392855447
$ws.Cells.Item(14, 11).Value = 1.067414476103898
$ws.Cells.Item(14, 12).Value = 1.070338416298749
$ws.Cells.Item(14, 13).Value = 1.07975045412107
$ws.Cells.Item(14, 14).Value = 1.069614213665785

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.061832223369714
$ws.Cells.Item(15, 4).Value = 1.064041207473477
$ws.Cells.Item(15, 5).Value = 1.066987557370471
$ws.Cells.Item(15, 6).Value = 1.076438585838582
$ws.Cells.Item(15, 9).Value = 1.050849440310343
$ws.Cells.Item(15, 10).Value = 1.068184540305053
$ws.Cells.Item(15, 11).Value = 1.067492074566953
$ws.Cells.Item(15, 12).Value = 1.070428206702333
$ws.Cells.Item(15, 13).Value = 1.079846874593457
$ws.Cells.Item(15, 14).Value = 1.069701484874766

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.062552248676424
$ws.Cells.Item(16, 4).Value = 1.064604954836024
$ws.Cells.Item(16, 5).Value = 1.067622275549515
$ws.Cells.Item(16, 6).Value = 1.077110928395892
$ws.Cells.Item(16, 9).Value = 1.051062566280474
$ws.Cells.Item(16, 10).Value = 1.068691609649342
$ws.Cells.Item(16, 11).Value = 1.067943531836963
$ws.Cells.Item(16, 12).Value = 1.070950756108582
$ws.Cells.Item(16, 13).Value = 1.080408062645753
$ws.Cells.Item(16, 14).Value = 1.070209274315685

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.063003940604004
$ws.Cells.Item(17, 4).Value = 1.064958612673301
$ws.Cells.Item(17, 5).Value = 1.068020557601454
$ws.Cells.Item(17, 6).Value = 1.077532850483434
$ws.Cells.Item(17, 9).Value = 1.051195885767508
$ws.Cells.Item(17, 10).Value = 1.069009534550417
$ws.Cells.Item(17, 11).Value = 1.068226543696577
$ws.Cells.Item(17, 12).Value = 1.071278476524045
$ws.Cells.Item(17, 13).Value = 1.080760063006025
$ws.Cells.Item(17, 14).Value = 1.070527650706586

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.063267419278342
$ws.Cells.Item(18, 4).Value = 1.065164907586602
$ws.Cells.Item(18, 5).Value = 1.0682529193474
$ws.Cells.Item(18, 6).Value = 1.077779015216554
$ws.Cells.Item(18, 9).Value = 1.051273515906358
$ws.Cells.Item(18, 10).Value = 1.069194922167942
$ws.Cells.Item(18, 11).Value = 1.06839155660585
$ws.Cells.Item(18, 12).Value = 1.071469608398147
$ws.Cells.Item(18, 13).Value = 1.080965372429097
$ws.Cells.Item(18, 14).Value = 1.070713301595792

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.063357261173949
$ws.Cells.Item(19, 4).Value = 1.065235250906041
$ws.Cells.Item(19, 5).Value = 1.068332157281145
$ws.Cells.Item(19, 6).Value = 1.077862962013681
$ws.Cells.Item(19, 9).Value = 1.051299963246474
$ws.Cells.Item(19, 10).Value = 1.069258125645084
$ws.Cells.Item(19, 11).Value = 1.068447811044555
$ws.Cells.Item(19, 12).Value = 1.071534775722051
$ws.Cells.Item(19, 13).Value = 1.081035376553684
$ws.Cells.Item(19, 14).Value = 1.070776594829121

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.062955476851225
$ws.Cells.Item(20, 4).Value = 1.064920667264552
$ws.Cells.Item(20, 5).Value = 1.067977820478033
$ws.Cells.Item(20, 6).Value = 1.077487575544011
$ws.Cells.Item(20, 9).Value = 1.051181595591889
$ws.Cells.Item(20, 10).Value = 1.068975429660566
$ws.Cells.Item(20, 11).Value = 1.068196185732871
$ws.Cells.Item(20, 12).Value = 1.071243317492103
$ws.Cells.Item(20, 13).Value = 1.080722297365681
$ws.Cells.Item(20, 14).Value = 1.07049349738388

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.061649404264505
$ws.Cells.Item(21, 4).Value = 1.063898069601242
$ws.Cells.Item(21, 5).Value = 1.066826431061867
$ws.Cells.Item(21, 6).Value = 1.076267917939926
$ws.Cells.Item(21, 9).Value = 1.050795208994763
$ws.Cells.Item(21, 10).Value = 1.06805573859314
$ws.Cells.Item(21, 11).Value = 1.06737738511328
$ws.Cells.Item(21, 12).Value = 1.070295500573795
$ws.Cells.Item(21, 13).Value = 1.079704370477644
$ws.Cells.Item(21, 14).Value = 1.069572500249647

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.060828347335986
$ws.Cells.Item(22, 4).Value = 1.063255231997178
$ws.Cells.Item(22, 5).Value = 1.066102960846398
$ws.Cells.Item(22, 6).Value = 1.075501652179592
$ws.Cells.Item(22, 9).Value = 1.050551074405516
$ws.Cells.Item(22, 10).Value = 1.067477015747591
$ws.Cells.Item(22, 11).Value = 1.066862002381864
$ws.Cells.Item(22, 12).Value = 1.069699371204555
$ws.Cells.Item(22, 13).Value = 1.079064299711193
$ws.Cells.Item(22, 14).Value = 1.068992955551278

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.061263593858578
$ws.Cells.Item(23, 4).Value = 1.063596001958772
$ws.Cells.Item(23, 5).Value = 1.066486443134703
$ws.Cells.Item(23, 6).Value = 1.075907809129143
$ws.Cells.Item(23, 9).Value = 1.050680608377407
$ws.Cells.Item(23, 10).Value = 1.067783852982256
$ws.Cells.Item(23, 11).Value = 1.067135270654552
$ws.Cells.Item(23, 12).Value = 1.070015409644183
$ws.Cells.Item(23, 13).Value = 1.079403619083273
$ws.Cells.Item(23, 14).Value = 1.06930022853001

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.062977375490128
$ws.Cells.Item(24, 4).Value = 1.06493781312334
$ws.Cells.Item(24, 5).Value = 1.067997131388102
$ws.Cells.Item(24, 6).Value = 1.077508033140149
$ws.Cells.Item(24, 9).Value = 1.051188053117563
$ws.Cells.Item(24, 10).Value = 1.068990840354586
$ws.Cells.Item(24, 11).Value = 1.068209903384799
$ws.Cells.Item(24, 12).Value = 1.071259204412726
$ws.Cells.Item(24, 13).Value = 1.080739362052646
$ws.Cells.Item(24, 14).Value = 1.070508929962852

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.064966692760567
$ws.Cells.Item(25, 4).Value = 1.066495392182126
$ws.Cells.Item(25, 5).Value = 1.069752185210875
$ws.Cells.Item(25, 6).Value = 1.079367543545352
$ws.Cells.Item(25, 9).Value = 1.051771703071795
$ws.Cells.Item(25, 10).Value = 1.070389420012588
$ws.Cells.Item(25, 11).Value = 1.069454479541153
$ws.Cells.Item(25, 12).Value = 1.072701700863774
$ws.Cells.Item(25, 13).Value = 1.082289165954997
$ws.Cells.Item(25, 14).Value = 1.071909495764389
